$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Scorpion"
$ws.Range("C2").Value = 10000

$ws.Range("C2").Select()
